$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header columns (values come from / extend sharedStrings.xml)
$ws.Range("I1").Value = "deposit"
$ws.Range("J1").Value = "balance"

# Match the header styling used by the rest of row 1 (blue fill, style index 1)
$ws.Range("I1:J1").Interior.Color = $ws.Range("H1").Interior.Color

# Update the active selection to match the committed state
[void]$ws.Range("K5").Select()
